$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.377.54"
$ws.Range('E2').Value = "  +1.94%  "
$ws.Range('D3').Value = "'3.391.87"
$ws.Range('E3').Value = "  +1.49%  "
$ws.Range('E4').Value = "  +0.03%  "
$ws.Range('D5').Value = "'588.36"
$ws.Range('E5').Value = "  +0.76%  "
$ws.Range('D6').Value = "'180.78"
$ws.Range('E6').Value = "  +1.81%  "
$ws.Range('E7').Value = "  +0.08%  "
$ws.Range('D8').Value = "'0.596"
$ws.Range('E8').Value = "  +0.69%  "
$ws.Range('E9').Value = "  +8.42%  "
$ws.Range('D10').Value = "'0.591"
$ws.Range('E10').Value = "  +1.57%  "
$ws.Range('D11').Value = "'48.76"
$ws.Range('E11').Value = "  +3.82%  "
$ws.Range('E12').Value = "  +3.59%  "
$ws.Range('D13').Value = "'683.94"
$ws.Range('E13').Value = "  -0.70%  "
$ws.Range('E14').Value = "  +2.06%  "
$ws.Range('D15').Value = "'3.936.73"
$ws.Range('E15').Value = "  +1.35%  "
$ws.Range('D16').Value = "'69.417.04"
$ws.Range('E16').Value = "  +2.01%  "
$ws.Range('D17').Value = "'3.419.99"
$ws.Range('E17').Value = "  +2.37%  "
$ws.Range('E18').Value = "  +1.75%  "
$ws.Range('D19').Value = "'17.73"
$ws.Range('E19').Value = "  +1.67%  "
$ws.Range('D20').Value = "'11.40"
$ws.Range('E20').Value = "  +2.64%  "
$ws.Range('D21').Value = "'0.903"
$ws.Range('E21').Value = "  +0.51%  "
$ws.Range('D22').Value = "'5.43"
$ws.Range('E22').Value = "  +1.23%  "
$ws.Range('E23').Value = "  +0.25%  "
$ws.Range('D24').Value = "'104.06"
$ws.Range('E24').Value = "  +5.58%  "
$ws.Range('D25').Value = "'3.93"
$ws.Range('E25').Value = "  +0.86%  "
$ws.Range('E26').Value = "  +1.30%  "
$ws.Range('E27').Value = "  +0.97%  "
$ws.Range('E28').Value = "  +3.43%  "
$ws.Range('E29').Value = "  +1.88%  "
$ws.Range('D30').Value = "'7.00"
$ws.Range('E30').Value = "  -1.57%  "
$ws.Range('E31').Value = "  +1.65%  "
$ws.Range('E32').Value = "  +9.37%  "
$ws.Range('D33').Value = "'555.12"
$ws.Range('E33').Value = "  -3.65%  "
$ws.Range('E34').Value = "  +0.75%  "
$ws.Range('D35').Value = "'58.16"
$ws.Range('E35').Value = "  +1.45%  "
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = "  -0.04%  "
$ws.Range('D37').Value = "'3.709.32"
$ws.Range('E37').Value = "  -0.44%  "
$ws.Range('E38').Value = "  +6.60%  "
$ws.Range('D39').Value = "'35.15"
$ws.Range('E39').Value = "  +1.82%  "
$ws.Range('E40').Value = "  +1.56%  "
$ws.Range('E41').Value = "  +4.11%  "
$ws.Range('E42').Value = "  +0.50%  "
$ws.Range('E43').Value = "  +0.97%  "
$ws.Range('D44').Value = "'0.0423"
$ws.Range('E44').Value = "  +4.09%  "
$ws.Range('E45').Value = "  -2.64%  "
$ws.Range('D46').Value = "'2.66"
$ws.Range('E46').Value = "  -0.21%  "
$ws.Range('E47').Value = "  +0.72%  "
$ws.Range('E48').Value = "  +4.75%  "
$ws.Range('E49').Value = "  +0.08%  "
$ws.Range('D50').Value = "'132.01"
$ws.Range('E50').Value = "  +2.02%  "
$ws.Range('E51').Value = "  -1.41%  "
